# "Check if valid template or not"
# Staff_Details -> Student_Details: rename headers and convert the
# Staff_Type (free-text) column into a numeric Student_Grade column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Student_Name"
$ws.Range("B1").Value = "Student_Id"
$ws.Range("C1").Value = "Student_Grade"
$ws.Range("D1").Value = "Student_Address"
$ws.Range("E1").Value = "Student_ZipCode"

# --- Column C becomes numeric (grade number instead of job title) ----
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 1

# --- Column widths (best effort - engine quantizes to 1/6 char units) -
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(4).ColumnWidth = 14
$ws.Columns.Item(5).ColumnWidth = 13.5

# --- Selection moved to F4 on the data sheet --------------------------
$ws.Range("F4").Select()
